$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the old "IDTYPE" column (column B) to hold step titles.
$ws.Columns.Item(2).Insert()

# Header + values for the new STEP_TITLE column.
$ws.Range("B1").Value = "STEP_TITLE"
$ws.Range("B2").Value = "Login"
$ws.Range("B3").Value = "Fill Username"
$ws.Range("B4").Value = "Fill Password"
$ws.Range("B5").Value = "Submit form"

# Recompute the (best-fit) column widths affected by the inserted column / now-empty columns.
$ws.Columns.Item(1).ColumnWidth = 10.333333333333332
$ws.Columns.Item(2).ColumnWidth = 11.166666666666666
$ws.Columns.Item(5).ColumnWidth = 5.666666666666667
$ws.Columns.Item(7).ColumnWidth = 8.5
$ws.Columns.Item(8).ColumnWidth = 10.166666666666666

# Match the author's last selected cell.
$ws.Range("D17").Select() | Out-Null
